# edit.ps1
# Reproduces the authored change to Blogs.xlsx:
#  - The "BlogPosts" sheet had its columns reordered so that "BlogName"
#    (previously column C) now comes first (column A), pushing
#    "Title"/"Slug" to column B and "Content" to column C.
#  - Row heights on BlogPosts rows 2-5 were tightened slightly (16.5 -> 16.2).
#  - Column widths on BlogPosts were adjusted to fit the new layout.
#  - The workbook was re-saved with the "BlogPosts" sheet active/selected
#    (instead of "Blogs"), with the last selected cell at C10.

$wb = $excel.ActiveWorkbook

$blogs = $wb.Worksheets.Item("Blogs")
$posts = $wb.Worksheets.Item("BlogPosts")

# --- BlogPosts: rotate the BlogName column (C) to the front (A) -----------
# Cutting column C and inserting it before column A shifts the old A->B and
# old B->C, giving: A=BlogName, B=Title/Slug, C=Content.
$posts.Columns.Item(3).Cut()
$posts.Columns.Item(1).Insert()

# --- BlogPosts: row heights for the data rows ------------------------------
$posts.Rows.Item(2).RowHeight = 16.2
$posts.Rows.Item(3).RowHeight = 16.2
$posts.Rows.Item(4).RowHeight = 16.2
$posts.Rows.Item(5).RowHeight = 16.2

# --- BlogPosts: column widths for the (now) Title and Content columns -----
$posts.Columns.Item(2).ColumnWidth = 45.142857142857146
$posts.Columns.Item(3).ColumnWidth = 48.857142857142854

# --- Make BlogPosts the active sheet / selection ---------------------------
$posts.Activate() | Out-Null
$posts.Range("C10").Select() | Out-Null
